# Fill in the placeholder "test" row (row 2) with proper test values and
# move the active cell selection, as per the commit
# "alles zu private, serialVersionId angepasst".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 currently has blanks in forename / surname / street_nr / city.
# Populate them with test data; keep username/email/password/role as-is.
$ws.Range("C2").Value = "TestVorname"
$ws.Range("D2").Value = "TestNachname"
$ws.Range("E2").Value = "Teststraße 0"
$ws.Range("G2").Value = "Teststadt"

# Re-write the numeric cells so they are stored as plain integers
# (0 / 0 / 1) instead of the old 0.0 / 0.0 / 1.0 representation.
$ws.Range("A2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("J2").Value = 1

# Move the active cell selection from M3 to M4.
$ws.Range("M4").Select()

$wb.Save()
